$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 79; this shifts the existing rows 79-102 down to 80-103
# and Excel extends the used range / dimension automatically.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new record. The surrounding
# columns (A,B,C,E,F,G,H,I,J,Q,T) keep the same values as the rest of this
# data block, only the fields below actually differ for this record.
$ws.Cells.Item(79, 1).Value = 5
$ws.Cells.Item(79, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(79, 3).Value = "Maule"
$ws.Cells.Item(79, 4).Value = Get-Date -Year 2022 -Month 3 -Day 16 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(79, 5).Value = 7
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100103
$ws.Cells.Item(79, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(79, 9).Value = 100103002
$ws.Cells.Item(79, 10).Value = "Ciruela"
$ws.Cells.Item(79, 11).Value = "Angeleno"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 230
$ws.Cells.Item(79, 14).Value = 8000
$ws.Cells.Item(79, 15).Value = 8000
$ws.Cells.Item(79, 16).Value = 8000
$ws.Cells.Item(79, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(79, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(79, 19).Value = 444
$ws.Cells.Item(79, 20).Value = 18

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(79, 4).NumberFormat = $ws.Cells.Item(80, 4).NumberFormat
